# Update 27/11/2023 - 5:19
# Rescale the stiffness-matrix-like values in Sheet1 (KGlob) that were
# previously entered/exported with an incorrect x1000 scale factor, and
# refresh a handful of dependent condensed entries accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.243503116200358
$ws.Range("K2").Value = -1.243503116200358
$ws.Range("C3").Value = 0.03454175322778771
$ws.Range("D3").Value = 0.05181262984168157
$ws.Range("L3").Value = -0.03454175322778771
$ws.Range("M3").Value = 0.05181262984168157
$ws.Range("C4").Value = 0.05181262984168157
$ws.Range("D4").Value = 0.1036252596833631
$ws.Range("L4").Value = -0.05181262984168157
$ws.Range("M4").Value = 0.05181262984168156
$ws.Range("E5").Value = 1.578881281869029
$ws.Range("H5").Value = -0.4898398753146541
$ws.Range("K5").Value = -0.5329299069430105
$ws.Range("N5").Value = 0
$ws.Range("S5").Value = -0.556111499611364
$ws.Range("F6").Value = 0.007919917241546575
$ws.Range("G6").Value = -0.007193970350291186
$ws.Range("I6").Value = -0.002111378772907992
$ws.Range("J6").Value = -0.008039890837502311
$ws.Range("L6").Value = -0.002719030137464339
$ws.Range("M6").Value = -0.009516605481125186
$ws.Range("O6").Value = -0.003089508331174243
$ws.Range("P6").Value = 0.01036252596833631
$ws.Range("F7").Value = -0.007193970350291186
$ws.Range("G7").Value = 0.1315734401557524
$ws.Range("I7").Value = 0.008039890837502311
$ws.Range("J7").Value = 0.02040999480477725
$ws.Range("L7").Value = 0.009516605481125186
$ws.Range("M7").Value = 0.0222054127892921
$ws.Range("O7").Value = -0.01036252596833631
$ws.Range("P7").Value = 0.02317131248380683
$ws.Range("E8").Value = -0.4898398753146541
$ws.Range("H8").Value = 1.535160290982592
$ws.Range("K8").Value = -0.4627127369362951
$ws.Range("Q8").Value = -0.5826076787316429
$ws.Range("F9").Value = -0.002111378772907992
$ws.Range("G9").Value = 0.008039890837502311
$ws.Range("I9").Value = 0.007443528991664547
$ws.Range("J9").Value = 0.0122393385094735
$ws.Range("L9").Value = -0.001779664372831904
$ws.Range("M9").Value = -0.00717405643961745
$ws.Range("N9").Value = 0.01137350411158864
$ws.Range("R9").Value = -0.003552485845924652
$ws.Range("S9").Value = 0
$ws.Range("F10").Value = -0.008039890837502311
$ws.Range("G10").Value = 0.02040999480477725
$ws.Range("I10").Value = 0.0122393385094735
$ws.Range("J10").Value = 0.1279300242485493
$ws.Range("L10").Value = 0.00717405643961745
$ws.Range("M10").Value = 0.01927969737234562
$ws.Range("N10").Value = 0.02427531994715179
$ws.Range("R10").Value = -0.01137350411158864
$ws.Range("S10").Value = 0
$ws.Range("B11").Value = -1.243503116200358
$ws.Range("E11").Value = -0.5329299069430105
$ws.Range("H11").Value = -0.4627127369362951
$ws.Range("K11").Value = 2.239145760079663
$ws.Range("C12").Value = -0.03454175322778771
$ws.Range("D12").Value = -0.05181262984168157
$ws.Range("F12").Value = -0.002719030137464339
$ws.Range("G12").Value = 0.009516605481125186
$ws.Range("I12").Value = -0.001779664372831904
$ws.Range("J12").Value = 0.00717405643961745
$ws.Range("L12").Value = 0.03904044773808395
$ws.Range("M12").Value = -0.03512196792093893
$ws.Range("C13").Value = 0.05181262984168157
$ws.Range("D13").Value = 0.05181262984168156
$ws.Range("F13").Value = -0.009516605481125186
$ws.Range("G13").Value = 0.0222054127892921
$ws.Range("I13").Value = -0.00717405643961745
$ws.Range("J13").Value = 0.01927969737234562
$ws.Range("L13").Value = -0.03512196792093893
$ws.Range("M13").Value = 0.1865954800066386
$ws.Range("E14").Value = 0
$ws.Range("I14").Value = 0.01137350411158864
$ws.Range("J14").Value = 0.02427531994715179
$ws.Range("N14").Value = 0.04855063989430358
$ws.Range("R14").Value = -0.01137350411158864
$ws.Range("F15").Value = -0.003089508331174243
$ws.Range("G15").Value = -0.01036252596833631
$ws.Range("O15").Value = 0.003089508331174243
$ws.Range("P15").Value = -0.01036252596833631
$ws.Range("F16").Value = 0.01036252596833631
$ws.Range("G16").Value = 0.02317131248380683
$ws.Range("O16").Value = -0.01036252596833631
$ws.Range("P16").Value = 0.04634262496761365
$ws.Range("H17").Value = -0.5826076787316429
$ws.Range("Q17").Value = 0.5826076787316429
$ws.Range("I18").Value = -0.003552485845924652
$ws.Range("J18").Value = -0.01137350411158864
$ws.Range("N18").Value = -0.01137350411158864
$ws.Range("R18").Value = 0.003552485845924652
$ws.Range("S18").Value = 0
$ws.Range("E19").Value = -0.556111499611364
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0.556111499611364
